$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 207, shifting the existing
# rows 207-308 down to 209-310.
$ws.Rows.Item(207).EntireRow.Insert()
$ws.Rows.Item(208).EntireRow.Insert()

# New row 207: Feria Lagunitas de Puerto Montt / Ciboulette entry dated 44960
$ws.Range("A207").Value = 4
$ws.Range("B207").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C207").Value = "Los Lagos"
$ws.Range("D207").Value = 44960
$ws.Range("E207").Value = 10
$ws.Range("F207").Value = 100112039
$ws.Range("G207").Value = "Ciboulette"
$ws.Range("H207").Value = "Sin especificar"
$ws.Range("I207").Value = "Primera"
$ws.Range("J207").Value = 120
$ws.Range("K207").Value = 8000
$ws.Range("L207").Value = 8000
$ws.Range("M207").Value = 8000
$ws.Range("N207").Value = "$/docena de atados"
$ws.Range("O207").Value = "Provincia de Cautín"
$ws.Range("P207").Value = 2667
$ws.Range("Q207").Value = 3
$ws.Range("R207").Value = "Hortaliza"

# New row 208: Feria Lagunitas de Puerto Montt / Ciboulette entry also dated 44960
$ws.Range("A208").Value = 4
$ws.Range("B208").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C208").Value = "Los Lagos"
$ws.Range("D208").Value = 44960
$ws.Range("E208").Value = 10
$ws.Range("F208").Value = 100112039
$ws.Range("G208").Value = "Ciboulette"
$ws.Range("H208").Value = "Sin especificar"
$ws.Range("I208").Value = "Primera"
$ws.Range("J208").Value = 240
$ws.Range("K208").Value = 3500
$ws.Range("L208").Value = 3500
$ws.Range("M208").Value = 3500
$ws.Range("N208").Value = "$/docena de atados"
$ws.Range("O208").Value = "Región Metropolitana"
$ws.Range("P208").Value = 1167
$ws.Range("Q208").Value = 3
$ws.Range("R208").Value = "Hortaliza"
